# Converting more to graph
# Adds two new reference sheets - "Cellular_regions" and "Objects_of_interest" -
# in between "Cell_types" and "cell_phenotype_type_categories".

$wb = $excel.ActiveWorkbook

# --- Create the two new worksheets in the correct position -----------------
$afterSheet = $wb.Worksheets.Item("Cell_types")
$regions = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$regions.Name = "Cellular_regions"

$objects = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $regions)
$objects.Name = "Objects_of_interest"

# --- Populate "Cellular_regions" (row-major, left to right) ----------------
$regions.Cells.Item(1,1).Value = "ID"
$regions.Cells.Item(1,2).Value = "Cellular_region"
$regions.Cells.Item(1,3).Value = "Ontological_identifier"
$regions.Cells.Item(1,4).Value = "Description"

$regions.Cells.Item(2,1).Value = 1
$regions.Cells.Item(2,2).Value = "Soma"
$regions.Cells.Item(2,3).Value = "NLX:154731"
$regions.Cells.Item(2,4).Value = "NULL"

$regions.Cells.Item(3,1).Value = 2
$regions.Cells.Item(3,2).Value = "Distal dendrite"
$regions.Cells.Item(3,3).Value = "NLX:154733"
$regions.Cells.Item(3,4).Value = "NULL"

$regions.Cells.Item(4,1).Value = 3
$regions.Cells.Item(4,2).Value = "Proximal dendrite"
$regions.Cells.Item(4,3).Value = "NLX:154734"
$regions.Cells.Item(4,4).Value = "NULL"

$regions.Cells.Item(5,1).Value = 4
$regions.Cells.Item(5,2).Value = "Dendrite"
$regions.Cells.Item(5,3).Value = "SAO:1211023249"
$regions.Cells.Item(5,4).Value = "NULL"

$regions.Cells.Item(6,1).Value = 5
$regions.Cells.Item(6,2).Value = "Dendritic spine"
$regions.Cells.Item(6,3).Value = "GO:0043197"
$regions.Cells.Item(6,4).Value = "NULL"

$regions.Cells.Item(7,1).Value = 6
$regions.Cells.Item(7,2).Value = "Dendritic shaft"
$regions.Cells.Item(7,3).Value = "SAO:2034472720"
$regions.Cells.Item(7,4).Value = "NULL"

$regions.Cells.Item(8,1).Value = 7
$regions.Cells.Item(8,2).Value = "Vesicle containing profile"
$regions.Cells.Item(8,3).Value = "NULL"
$regions.Cells.Item(8,4).Value = "NULL"

$regions.Cells.Item(9,1).Value = 8
$regions.Cells.Item(9,2).Value = "Dendritic region, unspecified"
$regions.Cells.Item(9,3).Value = "NULL"
$regions.Cells.Item(9,4).Value = "NULL"

$regions.Cells.Item(10,1).Value = 9
$regions.Cells.Item(10,2).Value = "Somatodendritic domain"
$regions.Cells.Item(10,3).Value = "NULL"
$regions.Cells.Item(10,4).Value = "NULL"

$regions.Cells.Item(11,1).Value = 10
$regions.Cells.Item(11,2).Value = "Neuron"
$regions.Cells.Item(11,3).Value = "NULL"
$regions.Cells.Item(11,4).Value = "NULL"

# Column widths for "Cellular_regions"
$regions.Columns.Item(2).ColumnWidth = 28.65
$regions.Columns.Item(3).ColumnWidth = 23.9

# Leftover selection on "Cellular_regions" (cursor left outside the table)
$regions.Range("G19").Select()

# --- Populate "Objects_of_interest" (row-major, left to right) -------------
$objects.Cells.Item(1,1).Value = "ID"
$objects.Cells.Item(1,2).Value = "Object_of_interest"
$objects.Cells.Item(1,3).Value = "Ontological_identifier"

$objects.Cells.Item(2,1).Value = 1
$objects.Cells.Item(2,2).Value = "Neurons"
$objects.Cells.Item(2,3).Value = "ILX:0107497"

$objects.Cells.Item(3,1).Value = 2
$objects.Cells.Item(3,2).Value = "Cells"
$objects.Cells.Item(3,3).Value = "ILX:0101839"

$objects.Cells.Item(4,1).Value = 4
$objects.Cells.Item(4,2).Value = "Axonal terminals"
$objects.Cells.Item(4,3).Value = "ILX:0101049"

$objects.Cells.Item(5,1).Value = 5
$objects.Cells.Item(5,2).Value = "Synapses"
$objects.Cells.Item(5,3).Value = "GO:0045202"

$objects.Cells.Item(6,1).Value = 9
$objects.Cells.Item(6,2).Value = "Synapses, symmetrical"
$objects.Cells.Item(6,3).Value = "ILX:0111392"

$objects.Cells.Item(7,1).Value = 10
$objects.Cells.Item(7,2).Value = "Synapses, asymmetrical"
$objects.Cells.Item(7,3).Value = "ILX:0100953"

$objects.Cells.Item(8,1).Value = 11
$objects.Cells.Item(8,2).Value = "Dendritic spines"
$objects.Cells.Item(8,3).Value = "GO:0043197"

$objects.Cells.Item(9,1).Value = 12
$objects.Cells.Item(9,2).Value = "Dendritic spines, mushroom"
$objects.Cells.Item(9,3).Value = "ILX:0107249"

$objects.Cells.Item(10,1).Value = 13
$objects.Cells.Item(10,2).Value = "Dendritic spines, stubby"
$objects.Cells.Item(10,3).Value = "ILX:0111129"

$objects.Cells.Item(11,1).Value = 14
$objects.Cells.Item(11,2).Value = "Dendritic spines, thin"
$objects.Cells.Item(11,3).Value = "ILX:0111691"

$objects.Cells.Item(12,1).Value = 17
$objects.Cells.Item(12,2).Value = "Dendritic spines, large"
$objects.Cells.Item(12,3).Value = "NULL"

$objects.Cells.Item(13,1).Value = 18
$objects.Cells.Item(13,2).Value = "Dendritic spines, giant"
$objects.Cells.Item(13,3).Value = "NULL"

$objects.Cells.Item(14,1).Value = 23
$objects.Cells.Item(14,2).Value = "Glia cell"
$objects.Cells.Item(14,3).Value = "NULL"

$objects.Cells.Item(15,1).Value = 24
$objects.Cells.Item(15,2).Value = "Axonal varicosities"
$objects.Cells.Item(15,3).Value = "NULL"

# Select the whole table and make this the active sheet/tab
$objects.Range("A1:C15").Select()
$objects.Activate()
